$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")
$ws.Range("B14").Value = "OLED screen from Ebay, USB powered fan from Amazon"
